$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 250001900
$ws.Range("J40").Value = 500000000
$ws.Range("L40").Value = 500000000
$ws.Range("N40").Value = -500000350

$ws.Range("H43").Value = 2384.25
$ws.Range("I43").Value = 2897.2
$ws.Range("J43").Value = 2017.8572
$ws.Range("K43").Value = 2897.2
$ws.Range("L43").Value = 2017.8572
$ws.Range("M43").Value = -2828.2
$ws.Range("N43").Value = -2155.8572

$ws.Range("H69").Value = 13526.333
$ws.Range("I69").Value = 6260
$ws.Range("K69").Value = 18780
$ws.Range("M69").Value = -17906

$ws.Range("H72").Value = 13526.333
$ws.Range("I72").Value = 6260
$ws.Range("K72").Value = 56340
$ws.Range("M72").Value = -51972

$ws.Range("H74").Value = 5469.6
$ws.Range("I74").Value = 5469.6
$ws.Range("K74").Value = 5469.6
$ws.Range("M74").Value = -4533.6

$ws.Range("H77").Value = 5469.6
$ws.Range("I77").Value = 5469.6
$ws.Range("K77").Value = 27348
$ws.Range("M77").Value = -22668

$ws.Range("H118").Value = 2351.4546
$ws.Range("I118").Value = 1509.4286
$ws.Range("J118").Value = 3825
$ws.Range("K118").Value = 4528.2858
$ws.Range("L118").Value = 11475
$ws.Range("M118").Value = -2871.2858
$ws.Range("N118").Value = -14789

$ws.Range("H135").Value = 1719.1852
$ws.Range("I135").Value = 583.5217
$ws.Range("K135").Value = 5251.6953
$ws.Range("M135").Value = -2716.6953

$ws.Range("H137").Value = 619514.75
$ws.Range("I137").Value = 839.0833
$ws.Range("K137").Value = 2517.2499
$ws.Range("M137").Value = 32.7501000000002

$ws.Range("H138").Value = 3145.0784
$ws.Range("J138").Value = 5830.4443
$ws.Range("L138").Value = 17491.3329
$ws.Range("N138").Value = -27771.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 6969
$ws.Range("I10").Value = 6969
$ws.Range("K10").Value = 6969
$ws.Range("M10").Value = -6799

$ws.Range("H45").Value = 4332868
$ws.Range("J45").Value = 5144
$ws.Range("L45").Value = 5144
$ws.Range("N45").Value = -5898

$ws.Range("H61").Value = 14004981
$ws.Range("I61").Value = 15390055
$ws.Range("K61").Value = 15390055
$ws.Range("M61").Value = -15389843

$ws.Range("H110").Value = 13124.25
$ws.Range("I110").Value = 16999.334
$ws.Range("J110").Value = 1499
$ws.Range("K110").Value = 16999.334
$ws.Range("L110").Value = 1499
$ws.Range("M110").Value = -14954.334
$ws.Range("N110").Value = -5589

$ws.Range("H122").Value = 2538.0435
$ws.Range("I122").Value = 1938.0714
$ws.Range("J122").Value = 3471.3333
$ws.Range("K122").Value = 5814.2142
$ws.Range("L122").Value = 10413.9999
$ws.Range("M122").Value = -3364.2142
$ws.Range("N122").Value = -15313.9999

$ws.Range("H132").Value = 3129698
$ws.Range("I132").Value = 4672.7856
$ws.Range("K132").Value = 14018.3568
$ws.Range("M132").Value = -11488.3568

$ws.Range("H136").Value = 14004981
$ws.Range("I136").Value = 15390055
$ws.Range("K136").Value = 46170165
$ws.Range("M136").Value = -46167615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 46684.75
$ws.Range("J81").Value = 46684.75
$ws.Range("L81").Value = 46684.75
$ws.Range("N81").Value = -48806.75

$ws.Range("H84").Value = 46684.75
$ws.Range("J84").Value = 46684.75
$ws.Range("L84").Value = 140054.25
$ws.Range("N84").Value = -150662.25

$ws.Range("H99").Value = 3200.9
$ws.Range("I99").Value = 2626.25
$ws.Range("J99").Value = 5499.5
$ws.Range("K99").Value = 2626.25
$ws.Range("L99").Value = 5499.5
$ws.Range("M99").Value = -1128.25
$ws.Range("N99").Value = -8495.5

$ws.Range("H107").Value = 4563.4
$ws.Range("I107").Value = 4828.25
$ws.Range("K107").Value = 4828.25
$ws.Range("M107").Value = -2908.25

$ws.Range("H134").Value = 3848250
$ws.Range("I134").Value = 1985.7727
$ws.Range("K134").Value = 5957.3181
$ws.Range("M134").Value = -3422.3181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24052414
$ws.Range("I31").Value = 34485050
$ws.Range("K31").Value = 34485050
$ws.Range("M31").Value = -34484755

$ws.Range("H34").Value = 24052414
$ws.Range("I34").Value = 34485050
$ws.Range("K34").Value = 34485050
$ws.Range("M34").Value = -34484848

$ws.Range("H58").Value = 2831.75
$ws.Range("I58").Value = 2816
$ws.Range("K58").Value = 2816
$ws.Range("M58").Value = -2613

$ws.Range("H99").Value = 12927.19
$ws.Range("I99").Value = 7479.615
$ws.Range("J99").Value = 21779.5
$ws.Range("K99").Value = 7479.615
$ws.Range("L99").Value = 21779.5
$ws.Range("M99").Value = -5981.615
$ws.Range("N99").Value = -24775.5

$ws.Range("H107").Value = 2329.5715
$ws.Range("I107").Value = 1051.1666
$ws.Range("K107").Value = 1051.1666
$ws.Range("M107").Value = 868.8334

$ws.Range("H126").Value = 12927.19
$ws.Range("I126").Value = 7479.615
$ws.Range("J126").Value = 21779.5
$ws.Range("K126").Value = 22438.845
$ws.Range("L126").Value = 65338.5
$ws.Range("M126").Value = -19968.845
$ws.Range("N126").Value = -70278.5

$ws.Range("H132").Value = 3355.05
$ws.Range("I132").Value = 3212.625
$ws.Range("K132").Value = 9637.875
$ws.Range("M132").Value = -7107.875

$ws.Range("H136").Value = 2831.75
$ws.Range("I136").Value = 2816
$ws.Range("K136").Value = 8448
$ws.Range("M136").Value = -5898

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4993
$ws.Range("J68").Value = 4993
$ws.Range("L68").Value = 14979
$ws.Range("N68").Value = -16601

$ws.Range("H71").Value = 4993
$ws.Range("J71").Value = 4993
$ws.Range("L71").Value = 44937
$ws.Range("N71").Value = -53049

$ws.Range("H92").Value = 266
$ws.Range("J92").Value = 266
$ws.Range("L92").Value = 798
$ws.Range("N92").Value = -3294

$ws.Range("H130").Value = 12371.363
$ws.Range("J130").Value = 14423.571
$ws.Range("L130").Value = 43270.713
$ws.Range("N130").Value = -53310.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10633
$ws.Range("J52").Value = 10633
$ws.Range("L52").Value = 10633
$ws.Range("N52").Value = -11151

$ws.Range("H122").Value = 5373.625
$ws.Range("I122").Value = 5427
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 16281
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -13831
$ws.Range("N122").Value = -19900

$ws.Range("H136").Value = 82246.5
$ws.Range("J136").Value = 82246.5
$ws.Range("L136").Value = 246739.5
$ws.Range("N136").Value = -251839.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4820.625
$ws.Range("I40").Value = 4833.5
$ws.Range("J40").Value = 4782
$ws.Range("K40").Value = 4833.5
$ws.Range("L40").Value = 4782
$ws.Range("M40").Value = -4697.5
$ws.Range("N40").Value = -5054

$ws.Range("H48").Value = 50022.5
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H61").Value = 76924970
$ws.Range("I61").Value = 83335340
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 83335340
$ws.Range("L61").Value = 400
$ws.Range("M61").Value = -83335138
$ws.Range("N61").Value = -804

$ws.Range("H93").Value = 2418887.8
$ws.Range("I93").Value = 1203.4736
$ws.Range("K93").Value = 1203.4736
$ws.Range("M93").Value = 44.52639999999997

$ws.Range("H113").Value = 76924970
$ws.Range("I113").Value = 83335340
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 83335340
$ws.Range("L113").Value = 400
$ws.Range("M113").Value = -83333170
$ws.Range("N113").Value = -4740

$ws.Range("H122").Value = 3260.1143
$ws.Range("I122").Value = 3165.8125
$ws.Range("J122").Value = 4266
$ws.Range("K122").Value = 9497.4375
$ws.Range("L122").Value = 12798
$ws.Range("M122").Value = -7047.4375
$ws.Range("N122").Value = -17698

$ws.Range("H132").Value = 6333.5
$ws.Range("I132").Value = 3333
$ws.Range("J132").Value = 7333.6665
$ws.Range("K132").Value = 9999
$ws.Range("L132").Value = 22000.9995
$ws.Range("M132").Value = -7469
$ws.Range("N132").Value = -27060.9995

$ws.Range("H136").Value = 4621
$ws.Range("I136").Value = 3367
$ws.Range("K136").Value = 10101
$ws.Range("M136").Value = -7551

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5705.9
$ws.Range("I122").Value = 5539.5713
$ws.Range("J122").Value = 6094
$ws.Range("K122").Value = 16618.7139
$ws.Range("L122").Value = 18282
$ws.Range("M122").Value = -14168.7139
$ws.Range("N122").Value = -23182

$ws.Range("H126").Value = 14788.9
$ws.Range("I126").Value = 14765.444
$ws.Range("K126").Value = 44296.33199999999
$ws.Range("M126").Value = -41826.33199999999

$ws.Range("H130").Value = 99333
$ws.Range("J130").Value = 99333
$ws.Range("L130").Value = 99333
$ws.Range("N130").Value = -109373

$ws.Range("H132").Value = 557462.8
$ws.Range("I132").Value = 2102.6155
$ws.Range("J132").Value = 2001399.4
$ws.Range("K132").Value = 6307.8465
$ws.Range("L132").Value = 6004198.199999999
$ws.Range("M132").Value = -3777.8465
$ws.Range("N132").Value = -6009258.199999999

$ws.Range("H136").Value = 325461.25
$ws.Range("I136").Value = 3010.3103
$ws.Range("J136").Value = 5001000
$ws.Range("K136").Value = 9030.930899999999
$ws.Range("L136").Value = 15003000
$ws.Range("M136").Value = -6480.930899999999
$ws.Range("N136").Value = -15008100
